$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Materials sheet: remove the old extra header row, then re-order rows ---
# Old layout:
#   row1: (blank), "Thermal Conductivity", "Thermal Diffusion", "Density"  (bold)
#   row2: "material", "k", "d", "rho"                                      (bold)
#   row3: filmgold
#   row4: gold
#   row5: pyrex
#   row6: quartz
# New layout:
#   row1: "material", "k", "d", "rho"  (bold)
#   row2: filmgold
#   row3: pyrex
#   row4: gold
#   row5: quartz
$ws.Rows.Item(1).Delete()

# Swap the (now) row3 (gold) and row4 (pyrex) so pyrex comes before gold.
$row3vals = @($ws.Cells.Item(3,1).Value2, $ws.Cells.Item(3,2).Value2, $ws.Cells.Item(3,3).Value2, $ws.Cells.Item(3,4).Value2)
$row4vals = @($ws.Cells.Item(4,1).Value2, $ws.Cells.Item(4,2).Value2, $ws.Cells.Item(4,3).Value2, $ws.Cells.Item(4,4).Value2)

for ($c = 1; $c -le 4; $c++) {
    $ws.Cells.Item(3, $c).Value2 = $row4vals[$c-1]
    $ws.Cells.Item(4, $c).Value2 = $row3vals[$c-1]
}

# --- Add the new Magnification sheet after Materials ---
$magSheet = $wb.Worksheets.Add([Type]::Missing, $ws)
$magSheet.Name = "Magnification"

$magSheet.Cells.Item(1,1).Value2 = "magnification"
$magSheet.Cells.Item(1,2).Value2 = "spotsize"

$magSheet.Cells.Item(2,1).Value2 = 50
$magSheet.Cells.Item(2,2).Value2 = 0.000002

$magSheet.Cells.Item(3,1).Value2 = 20
$magSheet.Cells.Item(3,2).Value2 = 0.000005

$magSheet.Cells.Item(4,1).Value2 = 10
$magSheet.Cells.Item(4,2).Value2 = 0.00001

$magSheet.Cells.Item(5,1).Value2 = 5
$magSheet.Cells.Item(5,2).Value2 = 0.00002

$magSheet.Range("B2:B5").NumberFormat = "0.00E+00"

[void]$magSheet.Range("F12").Select()

# --- Fix the frozen pane / selection on the Materials sheet, and make it the active tab ---
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
[void]$ws.Range("A2").Select()
$win.FreezePanes = $true
[void]$ws.Rows.Item(3).Select()
